$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new column D values ("U") for rows 2-5
$ws.Range("D2").Value = "U"
$ws.Range("D3").Value = "U"
$ws.Range("D4").Value = "U"
$ws.Range("D5").Value = "U"

# Update the selection to G10
$ws.Range("G10").Select()
